# Add 2022-Q3 data
# 1) Insert a new row into the "总计" (summary) sheet for the 2022-Q3 period,
#    pushing the existing quarter rows down by one.
# 2) Insert a brand new "2022-Q3" worksheet (positioned right before "2022-Q2")
#    holding the underlying fund-holding detail for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert the 2022-Q3 row at row 2
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# New row picks up some stray formatting from the insert - strip it back off
# for the non-index columns (they are unstyled in the rest of the sheet).
$summary.Range("B2:D2").ClearFormats()

# Column A carries the bold/border/centered style used throughout the sheet;
# clone it from the row below (which held that style before the insert).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# Part 2: brand-new "2022-Q3" worksheet with the quarter's fund holdings
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Single data row
$q3.Range("A2").Value = 0
$q3.Range("A2").Font.Bold = $true
$q3.Range("A2").HorizontalAlignment = -4108
$q3.Range("A2").VerticalAlignment = -4160

# Text-like numeric-looking fields keep their original string formatting
# (leading zeros / trailing zeros must survive), so force text before writing.
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "002020"

$q3.Range("C2").Value = "国都创新驱动灵活配置混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.12"

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "65.45"

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.31"

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0040"

$q3.Range("H2").Value = 7
